$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-column style (style index used by D2:D81) onto the two brand-new rows
$dateFormat = $ws.Range("D81").NumberFormat
$ws.Range("D82").NumberFormat = $dateFormat
$ws.Range("D83").NumberFormat = $dateFormat

# Row 57
$ws.Cells.Item(57, 4).Value = 45006
$ws.Cells.Item(57, 12).Value = 'Primera'
$ws.Cells.Item(57, 13).Value = 80
$ws.Cells.Item(57, 14).Value = 18000
$ws.Cells.Item(57, 15).Value = 18000
$ws.Cells.Item(57, 16).Value = 18000
$ws.Cells.Item(57, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(57, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(57, 19).Value = 1125
$ws.Cells.Item(57, 20).Value = 16

# Row 58
$ws.Cells.Item(58, 4).Value = 44294
$ws.Cells.Item(58, 12).Value = 'Primera'
$ws.Cells.Item(58, 13).Value = 80
$ws.Cells.Item(58, 14).Value = 14000
$ws.Cells.Item(58, 15).Value = 16000
$ws.Cells.Item(58, 16).Value = 15000
$ws.Cells.Item(58, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(58, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(58, 19).Value = 938
$ws.Cells.Item(58, 20).Value = 16

# Row 59
$ws.Cells.Item(59, 4).Value = 44637
$ws.Cells.Item(59, 12).Value = 'Primera'
$ws.Cells.Item(59, 13).Value = 65
$ws.Cells.Item(59, 14).Value = 15000
$ws.Cells.Item(59, 15).Value = 15000
$ws.Cells.Item(59, 16).Value = 15000
$ws.Cells.Item(59, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(59, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(59, 19).Value = 938
$ws.Cells.Item(59, 20).Value = 16

# Row 60
$ws.Cells.Item(60, 4).Value = 44244
$ws.Cells.Item(60, 12).Value = 'Primera'
$ws.Cells.Item(60, 13).Value = 25
$ws.Cells.Item(60, 14).Value = 14000
$ws.Cells.Item(60, 15).Value = 14000
$ws.Cells.Item(60, 16).Value = 14000
$ws.Cells.Item(60, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(60, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(60, 19).Value = 875
$ws.Cells.Item(60, 20).Value = 16

# Row 61
$ws.Cells.Item(61, 4).Value = 44258
$ws.Cells.Item(61, 12).Value = 'Primera'
$ws.Cells.Item(61, 13).Value = 65
$ws.Cells.Item(61, 14).Value = 16000
$ws.Cells.Item(61, 15).Value = 16000
$ws.Cells.Item(61, 16).Value = 16000
$ws.Cells.Item(61, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(61, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(61, 19).Value = 889
$ws.Cells.Item(61, 20).Value = 18

# Row 62
$ws.Cells.Item(62, 4).Value = 44658
$ws.Cells.Item(62, 12).Value = 'Primera'
$ws.Cells.Item(62, 13).Value = 200
$ws.Cells.Item(62, 14).Value = 17000
$ws.Cells.Item(62, 15).Value = 17000
$ws.Cells.Item(62, 16).Value = 17000
$ws.Cells.Item(62, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(62, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(62, 19).Value = 1062
$ws.Cells.Item(62, 20).Value = 16

# Row 63
$ws.Cells.Item(63, 4).Value = 44273
$ws.Cells.Item(63, 12).Value = 'Primera'
$ws.Cells.Item(63, 13).Value = 55
$ws.Cells.Item(63, 14).Value = 14000
$ws.Cells.Item(63, 15).Value = 14000
$ws.Cells.Item(63, 16).Value = 14000
$ws.Cells.Item(63, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(63, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(63, 19).Value = 875
$ws.Cells.Item(63, 20).Value = 16

# Row 64
$ws.Cells.Item(64, 4).Value = 44389
$ws.Cells.Item(64, 12).Value = 'Primera'
$ws.Cells.Item(64, 13).Value = 30
$ws.Cells.Item(64, 14).Value = 30000
$ws.Cells.Item(64, 15).Value = 30000
$ws.Cells.Item(64, 16).Value = 30000
$ws.Cells.Item(64, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(64, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(64, 19).Value = 1875
$ws.Cells.Item(64, 20).Value = 16

# Row 65
$ws.Cells.Item(65, 4).Value = 44832
$ws.Cells.Item(65, 12).Value = 'Especial'
$ws.Cells.Item(65, 13).Value = 50
$ws.Cells.Item(65, 14).Value = 40000
$ws.Cells.Item(65, 15).Value = 40000
$ws.Cells.Item(65, 16).Value = 40000
$ws.Cells.Item(65, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(65, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(65, 19).Value = 2222
$ws.Cells.Item(65, 20).Value = 16

# Row 66
$ws.Cells.Item(66, 4).Value = 44435
$ws.Cells.Item(66, 12).Value = 'Especial'
$ws.Cells.Item(66, 13).Value = 30
$ws.Cells.Item(66, 14).Value = 45000
$ws.Cells.Item(66, 15).Value = 45000
$ws.Cells.Item(66, 16).Value = 45000
$ws.Cells.Item(66, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(66, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(66, 19).Value = 2500
$ws.Cells.Item(66, 20).Value = 18

# Row 67
$ws.Cells.Item(67, 4).Value = 44425
$ws.Cells.Item(67, 12).Value = 'Especial'
$ws.Cells.Item(67, 13).Value = 35
$ws.Cells.Item(67, 14).Value = 40000
$ws.Cells.Item(67, 15).Value = 40000
$ws.Cells.Item(67, 16).Value = 40000
$ws.Cells.Item(67, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(67, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(67, 19).Value = 2500
$ws.Cells.Item(67, 20).Value = 18

# Row 68
$ws.Cells.Item(68, 4).Value = 44987
$ws.Cells.Item(68, 12).Value = 'Primera'
$ws.Cells.Item(68, 13).Value = 40
$ws.Cells.Item(68, 14).Value = 15000
$ws.Cells.Item(68, 15).Value = 15000
$ws.Cells.Item(68, 16).Value = 15000
$ws.Cells.Item(68, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(68, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(68, 19).Value = 938
$ws.Cells.Item(68, 20).Value = 16

# Row 69
$ws.Cells.Item(69, 4).Value = 45005
$ws.Cells.Item(69, 12).Value = 'Primera'
$ws.Cells.Item(69, 13).Value = 100
$ws.Cells.Item(69, 14).Value = 18000
$ws.Cells.Item(69, 15).Value = 18000
$ws.Cells.Item(69, 16).Value = 18000
$ws.Cells.Item(69, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(69, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(69, 19).Value = 1125
$ws.Cells.Item(69, 20).Value = 16

# Row 70
$ws.Cells.Item(70, 4).Value = 44603
$ws.Cells.Item(70, 12).Value = 'Primera'
$ws.Cells.Item(70, 13).Value = 40
$ws.Cells.Item(70, 14).Value = 17000
$ws.Cells.Item(70, 15).Value = 18000
$ws.Cells.Item(70, 16).Value = 17500
$ws.Cells.Item(70, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(70, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(70, 19).Value = 1094
$ws.Cells.Item(70, 20).Value = 16

# Row 71
$ws.Cells.Item(71, 4).Value = 44650
$ws.Cells.Item(71, 12).Value = 'Primera'
$ws.Cells.Item(71, 13).Value = 45
$ws.Cells.Item(71, 14).Value = 18000
$ws.Cells.Item(71, 15).Value = 20000
$ws.Cells.Item(71, 16).Value = 19333
$ws.Cells.Item(71, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(71, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(71, 19).Value = 1208
$ws.Cells.Item(71, 20).Value = 16

# Row 72
$ws.Cells.Item(72, 4).Value = 44238
$ws.Cells.Item(72, 12).Value = 'Primera'
$ws.Cells.Item(72, 13).Value = 65
$ws.Cells.Item(72, 14).Value = 14000
$ws.Cells.Item(72, 15).Value = 14000
$ws.Cells.Item(72, 16).Value = 14000
$ws.Cells.Item(72, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(72, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(72, 19).Value = 875
$ws.Cells.Item(72, 20).Value = 16

# Row 73
$ws.Cells.Item(73, 4).Value = 44663
$ws.Cells.Item(73, 12).Value = 'Primera'
$ws.Cells.Item(73, 13).Value = 55
$ws.Cells.Item(73, 14).Value = 16000
$ws.Cells.Item(73, 15).Value = 16000
$ws.Cells.Item(73, 16).Value = 16000
$ws.Cells.Item(73, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(73, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(73, 19).Value = 1000
$ws.Cells.Item(73, 20).Value = 16

# Row 74
$ws.Cells.Item(74, 4).Value = 44868
$ws.Cells.Item(74, 12).Value = 'Primera'
$ws.Cells.Item(74, 13).Value = 30
$ws.Cells.Item(74, 14).Value = 40000
$ws.Cells.Item(74, 15).Value = 40000
$ws.Cells.Item(74, 16).Value = 40000
$ws.Cells.Item(74, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(74, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(74, 19).Value = 2500
$ws.Cells.Item(74, 20).Value = 16

# Row 75
$ws.Cells.Item(75, 4).Value = 44253
$ws.Cells.Item(75, 12).Value = 'Primera'
$ws.Cells.Item(75, 13).Value = 55
$ws.Cells.Item(75, 14).Value = 16000
$ws.Cells.Item(75, 15).Value = 16000
$ws.Cells.Item(75, 16).Value = 16000
$ws.Cells.Item(75, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(75, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(75, 19).Value = 1000
$ws.Cells.Item(75, 20).Value = 16

# Row 76
$ws.Cells.Item(76, 4).Value = 44998
$ws.Cells.Item(76, 12).Value = 'Primera'
$ws.Cells.Item(76, 13).Value = 35
$ws.Cells.Item(76, 14).Value = 15000
$ws.Cells.Item(76, 15).Value = 15000
$ws.Cells.Item(76, 16).Value = 15000
$ws.Cells.Item(76, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(76, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(76, 19).Value = 938
$ws.Cells.Item(76, 20).Value = 16

# Row 77
$ws.Cells.Item(77, 4).Value = 44957
$ws.Cells.Item(77, 12).Value = 'Primera'
$ws.Cells.Item(77, 13).Value = 25
$ws.Cells.Item(77, 14).Value = 30000
$ws.Cells.Item(77, 15).Value = 30000
$ws.Cells.Item(77, 16).Value = 30000
$ws.Cells.Item(77, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(77, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(77, 19).Value = 1667
$ws.Cells.Item(77, 20).Value = 18

# Row 78
$ws.Cells.Item(78, 4).Value = 44306
$ws.Cells.Item(78, 12).Value = 'Especial'
$ws.Cells.Item(78, 13).Value = 50
$ws.Cells.Item(78, 14).Value = 22000
$ws.Cells.Item(78, 15).Value = 22000
$ws.Cells.Item(78, 16).Value = 22000
$ws.Cells.Item(78, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(78, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(78, 19).Value = 1375
$ws.Cells.Item(78, 20).Value = 16

# Row 79
$ws.Cells.Item(79, 4).Value = 44606
$ws.Cells.Item(79, 12).Value = 'Segunda'
$ws.Cells.Item(79, 13).Value = 80
$ws.Cells.Item(79, 14).Value = 10000
$ws.Cells.Item(79, 15).Value = 10000
$ws.Cells.Item(79, 16).Value = 10000
$ws.Cells.Item(79, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(79, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(79, 19).Value = 625
$ws.Cells.Item(79, 20).Value = 18

# Row 80
$ws.Cells.Item(80, 4).Value = 44677
$ws.Cells.Item(80, 12).Value = 'Primera'
$ws.Cells.Item(80, 13).Value = 55
$ws.Cells.Item(80, 14).Value = 20000
$ws.Cells.Item(80, 15).Value = 20000
$ws.Cells.Item(80, 16).Value = 20000
$ws.Cells.Item(80, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(80, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(80, 19).Value = 1250
$ws.Cells.Item(80, 20).Value = 16

# Row 81
$ws.Cells.Item(81, 4).Value = 44236
$ws.Cells.Item(81, 12).Value = 'Primera'
$ws.Cells.Item(81, 13).Value = 45
$ws.Cells.Item(81, 14).Value = 25000
$ws.Cells.Item(81, 15).Value = 25000
$ws.Cells.Item(81, 16).Value = 25000
$ws.Cells.Item(81, 17).Value = '$/caja 16 kilos'
$ws.Cells.Item(81, 18).Value = 'Provincia de Limarí'
$ws.Cells.Item(81, 19).Value = 1562
$ws.Cells.Item(81, 20).Value = 16

# Row 82
$ws.Cells.Item(82, 1).Value = 10
$ws.Cells.Item(82, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(82, 3).Value = 'La Araucanía'
$ws.Cells.Item(82, 5).Value = 9
$ws.Cells.Item(82, 6).Value = 'Fruta'
$ws.Cells.Item(82, 7).Value = 100107
$ws.Cells.Item(82, 8).Value = 'Otros'
$ws.Cells.Item(82, 9).Value = 100107011
$ws.Cells.Item(82, 10).Value = 'Tuna'
$ws.Cells.Item(82, 11).Value = 'Sin especificar'
$ws.Cells.Item(82, 4).Value = 44601
$ws.Cells.Item(82, 12).Value = 'Especial'
$ws.Cells.Item(82, 13).Value = 30
$ws.Cells.Item(82, 14).Value = 25000
$ws.Cells.Item(82, 15).Value = 25000
$ws.Cells.Item(82, 16).Value = 25000
$ws.Cells.Item(82, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(82, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(82, 19).Value = 1389
$ws.Cells.Item(82, 20).Value = 18

# Row 83
$ws.Cells.Item(83, 1).Value = 10
$ws.Cells.Item(83, 2).Value = 'Vega Modelo de Temuco'
$ws.Cells.Item(83, 3).Value = 'La Araucanía'
$ws.Cells.Item(83, 5).Value = 9
$ws.Cells.Item(83, 6).Value = 'Fruta'
$ws.Cells.Item(83, 7).Value = 100107
$ws.Cells.Item(83, 8).Value = 'Otros'
$ws.Cells.Item(83, 9).Value = 100107011
$ws.Cells.Item(83, 10).Value = 'Tuna'
$ws.Cells.Item(83, 11).Value = 'Sin especificar'
$ws.Cells.Item(83, 4).Value = 44601
$ws.Cells.Item(83, 12).Value = 'Primera'
$ws.Cells.Item(83, 13).Value = 80
$ws.Cells.Item(83, 14).Value = 18000
$ws.Cells.Item(83, 15).Value = 18000
$ws.Cells.Item(83, 16).Value = 18000
$ws.Cells.Item(83, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(83, 18).Value = 'Provincia de Los Andes'
$ws.Cells.Item(83, 19).Value = 1000
$ws.Cells.Item(83, 20).Value = 18

